$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 41 - Cash And Cash Equivalents
$ws.Range("D41").Value = 20800

# Row 43 - Net Receivables
$ws.Range("I43").Value = 64100
$ws.Range("J43").Value = 56800

# Row 45 - Other Current Assets
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "NA"
$ws.Range("H45").Value = "NA"
$ws.Range("I45").Value = 7300
$ws.Range("J45").Value = 2000

# Row 47 - Total Current Assets
$ws.Range("D47").Value = 37400
$ws.Range("I47").Value = 129700

# Row 48 - Long Term Investments
$ws.Range("D48").Value = 229800

# Row 49 - Property Plant and Equipment
$ws.Range("D49").Value = 247800

# Row 52 - Deferred Long Term Asset Charges
$ws.Range("D52").Value = 136500

# Row 57 - Accounts Payable
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0

# Row 58 - Short/Current Long Term Debt
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0

# Row 59 - Other Current Liabilities
$ws.Range("D59").Value = 50300
$ws.Range("I59").Value = 87600
$ws.Range("J59").Value = 67600

# Row 61 - Long Term Debt
$ws.Range("J61").Value = 56000

# Row 62 - Other Liabilities
$ws.Range("D62").Value = 28700
$ws.Range("J62").Value = 2700
